# Apply 2025-07-03 data update to violent-crime-full-year workbook
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item('Citywide Totals')
$ws.Range('K2').Value = 7901
$ws.Range('L2').Value = 3302
$ws.Range('L3').Value = 3419
$ws.Range('I4').Value = 1838
$ws.Range('L4').Value = 857
$ws.Range('L5').Value = 195
$ws.Range('L6').Value = 3017
$ws.Range('I7').Value = 26306
$ws.Range('K7').Value = 27557
$ws.Range('L7').Value = 10790

$ws = $wb.Worksheets.Item('By Neighborhood')
$ws.Range('L7').Value = 363
$ws.Range('L8').Value = 686
$ws.Range('L9').Value = 69
$ws.Range('L11').Value = 181
$ws.Range('L12').Value = 24
$ws.Range('L14').Value = 53
$ws.Range('L15').Value = 81
$ws.Range('L19').Value = 304
$ws.Range('L23').Value = 118
$ws.Range('L27').Value = 102
$ws.Range('L29').Value = 582
$ws.Range('L30').Value = 57
$ws.Range('L31').Value = 102
$ws.Range('L33').Value = 499
$ws.Range('L35').Value = 17
$ws.Range('L36').Value = 145
$ws.Range('L37').Value = 387
$ws.Range('L41').Value = 47
$ws.Range('L42').Value = 345
$ws.Range('L47').Value = 78
$ws.Range('L48').Value = 154
$ws.Range('L51').Value = 132
$ws.Range('L54').Value = 226
$ws.Range('L57').Value = 39
$ws.Range('L60').Value = 65
$ws.Range('I63').Value = 262
$ws.Range('L63').Value = 37
$ws.Range('L65').Value = 209
$ws.Range('K67').Value = 1071
$ws.Range('L73').Value = 94
$ws.Range('L74').Value = 10
$ws.Range('L76').Value = 153
$ws.Range('L78').Value = 136
$ws.Range('L79').Value = 277
$ws.Range('L83').Value = 254
$ws.Range('L84').Value = 105
$ws.Range('L85').Value = 547
$ws.Range('L86').Value = 82
$ws.Range('L91').Value = 153
$ws.Range('L94').Value = 128
$ws.Range('L95').Value = 149
$ws.Range('L96').Value = 104
$ws.Range('L97').Value = 94
$ws.Range('I101').Value = 26306
$ws.Range('K101').Value = 27557
$ws.Range('L101').Value = 10790

$ws = $wb.Worksheets.Item('Bridgeport')
$ws.Range('L6').Value = 12
$ws.Range('L7').Value = 53

$ws = $wb.Worksheets.Item('West Ridge')
$ws.Range('L2').Value = 39
$ws.Range('L7').Value = 104

$ws = $wb.Worksheets.Item('Auburn Gresham')
$ws.Range('L2').Value = 115
$ws.Range('L3').Value = 112
$ws.Range('L6').Value = 103
$ws.Range('L7').Value = 363

$ws = $wb.Worksheets.Item('Belmont Cragin')
$ws.Range('L3').Value = 54
$ws.Range('L4').Value = 14
$ws.Range('L6').Value = 43
$ws.Range('L7').Value = 181

$ws = $wb.Worksheets.Item('South Shore')
$ws.Range('L3').Value = 223
$ws.Range('L4').Value = 40
$ws.Range('L6').Value = 116
$ws.Range('L7').Value = 547

$ws = $wb.Worksheets.Item('Little Village')
$ws.Range('L2').Value = 76
$ws.Range('L6').Value = 54

$ws = $wb.Worksheets.Item('Austin')
$ws.Range('L2').Value = 202
$ws.Range('L3').Value = 224
$ws.Range('L6').Value = 192
$ws.Range('L7').Value = 686

$ws = $wb.Worksheets.Item('South Chicago')
$ws.Range('L2').Value = 84
$ws.Range('L3').Value = 101
$ws.Range('L7').Value = 254

$ws = $wb.Worksheets.Item('Garfield Park')
$ws.Range('L2').Value = 143
$ws.Range('L3').Value = 151
$ws.Range('L6').Value = 172
$ws.Range('L7').Value = 499

$ws = $wb.Worksheets.Item('West Pullman')
$ws.Range('L5').Value = 6
$ws.Range('L7').Value = 149

$ws = $wb.Worksheets.Item('Grand Crossing')
$ws.Range('L4').Value = 24
$ws.Range('L7').Value = 387

$ws = $wb.Worksheets.Item('New City')
$ws.Range('L3').Value = 63
$ws.Range('L7').Value = 209

$ws = $wb.Worksheets.Item('Fuller Park')
$ws.Range('L3').Value = 11
$ws.Range('L7').Value = 57

$ws = $wb.Worksheets.Item('Gage Park')
$ws.Range('L6').Value = 34
$ws.Range('L7').Value = 102

$ws = $wb.Worksheets.Item('North Lawndale')
$ws.Range('K2').Value = 299
$ws.Range('K7').Value = 1071

$ws = $wb.Worksheets.Item('South Deering')
$ws.Range('L6').Value = 24
$ws.Range('L7').Value = 105

$ws = $wb.Worksheets.Item('Loop')
$ws.Range('L2').Value = 48
$ws.Range('L3').Value = 52
$ws.Range('L4').Value = 18
$ws.Range('L7').Value = 226

$ws = $wb.Worksheets.Item('Englewood')
$ws.Range('L2').Value = 174
$ws.Range('L3').Value = 221
$ws.Range('L4').Value = 28
$ws.Range('L6').Value = 151
$ws.Range('L7').Value = 582

$ws = $wb.Worksheets.Item('Lake View')
$ws.Range('L3').Value = 35
$ws.Range('L6').Value = 66
$ws.Range('L7').Value = 154

$ws = $wb.Worksheets.Item('Chatham')
$ws.Range('L3').Value = 92
$ws.Range('L5').Value = 5
$ws.Range('L7').Value = 304

$ws = $wb.Worksheets.Item('River North')
$ws.Range('L6').Value = 69
$ws.Range('L7').Value = 153

$ws = $wb.Worksheets.Item('Hermosa')
$ws.Range('L3').Value = 18
$ws.Range('L7').Value = 47

$ws = $wb.Worksheets.Item('Humboldt Park')
$ws.Range('L3').Value = 110
$ws.Range('L7').Value = 345

$ws = $wb.Worksheets.Item('Rogers Park')
$ws.Range('L2').Value = 39
$ws.Range('L7').Value = 136

$ws = $wb.Worksheets.Item('Douglas')
$ws.Range('L6').Value = 29
$ws.Range('L7').Value = 118

$ws = $wb.Worksheets.Item('Washington Park')
$ws.Range('L6').Value = 22
$ws.Range('L7').Value = 153

$ws = $wb.Worksheets.Item('Roseland')
$ws.Range('L3').Value = 100
$ws.Range('L7').Value = 277

$ws = $wb.Worksheets.Item('Grand Boulevard')
$ws.Range('L3').Value = 42
$ws.Range('L6').Value = 38
$ws.Range('L7').Value = 145

$ws = $wb.Worksheets.Item('West Loop')
$ws.Range('L6').Value = 47
$ws.Range('L7').Value = 128

$ws = $wb.Worksheets.Item('Kenwood')
$ws.Range('L6').Value = 18
$ws.Range('L7').Value = 78

$ws = $wb.Worksheets.Item('Brighton Park')
$ws.Range('L2').Value = 33
$ws.Range('L7').Value = 81

$ws = $wb.Worksheets.Item('Gold Coast')
$ws.Range('L6').Value = 7
$ws.Range('L7').Value = 17

$ws = $wb.Worksheets.Item('Avalon Park')
$ws.Range('L2').Value = 20
$ws.Range('L7').Value = 69

$ws = $wb.Worksheets.Item('Portage Park')
$ws.Range('L2').Value = 37
$ws.Range('L6').Value = 24
$ws.Range('L7').Value = 94

$ws = $wb.Worksheets.Item('West Town')
$ws.Range('L6').Value = 52
$ws.Range('L7').Value = 94

$ws = $wb.Worksheets.Item('Edgewater')
$ws.Range('L2').Value = 25
$ws.Range('L7').Value = 102

$ws = $wb.Worksheets.Item('Streeterville')
$ws.Range('L4').Value = 44
$ws.Range('L7').Value = 82

$ws = $wb.Worksheets.Item('Little Italy, UIC')
$ws.Range('L2').Value = 39
$ws.Range('L7').Value = 132

$ws = $wb.Worksheets.Item('Mckinley Park')
$ws.Range('L3').Value = 11
$ws.Range('L7').Value = 39

$ws = $wb.Worksheets.Item('Morgan Park')
$ws.Range('L3').Value = 26
$ws.Range('L7').Value = 65

$ws = $wb.Worksheets.Item('Beverly')
$ws.Range('L4').Value = 5
$ws.Range('L7').Value = 24

$ws = $wb.Worksheets.Item('Printers Row')
$ws.Range('L3').Value = 2
$ws.Range('L7').Value = 10
